$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.764.52'
$ws.Range("E2").Value = '  -4.98%  '
$ws.Range("D3").Value = '3.169.24'
$ws.Range("E3").Value = '  -5.30%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '531.62'
$ws.Range("E5").Value = '  -6.17%  '
$ws.Range("D6").Value = '134.61'
$ws.Range("E6").Value = '  -7.99%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.169.72'
$ws.Range("E8").Value = '  -5.32%  '
$ws.Range("D9").Value = '0.453'
$ws.Range("E9").Value = '  -6.33%  '
$ws.Range("D10").Value = '7.26'
$ws.Range("E10").Value = '  -8.22%  '
$ws.Range("E11").Value = '  -8.55%  '
$ws.Range("D12").Value = '0.395'
$ws.Range("E12").Value = '  -4.70%  '
$ws.Range("D13").Value = '3.715.40'
$ws.Range("E13").Value = '  -5.15%  '
$ws.Range("D15").Value = '25.78'
$ws.Range("E15").Value = '  -6.62%  '
$ws.Range("D16").Value = '3.175.03'
$ws.Range("E16").Value = '  -5.20%  '
$ws.Range("D17").Value = '58.139.09'
$ws.Range("E17").Value = '  -4.32%  '
$ws.Range("D18").Value = '0.0000155'
$ws.Range("E18").Value = '  -8.37%  '
$ws.Range("D19").Value = '5.83'
$ws.Range("E19").Value = '  -6.96%  '
$ws.Range("D20").Value = '13.25'
$ws.Range("E20").Value = '  -8.33%  '
$ws.Range("D21").Value = '8.08'
$ws.Range("E21").Value = '  -8.94%  '
$ws.Range("D22").Value = '357.73'
$ws.Range("E22").Value = '  -4.75%  '
$ws.Range("D23").Value = '0.997'
$ws.Range("E23").Value = '  -0.42%  '
$ws.Range("D24").Value = '69.57'
$ws.Range("E24").Value = '  -6.99%  '
$ws.Range("D25").Value = '0.516'
$ws.Range("E25").Value = '  -7.72%  '
$ws.Range("D26").Value = '3.314.07'
$ws.Range("E26").Value = '  -5.22%  '
$ws.Range("E27").Value = '  -3.52%  '
$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").Value = '0.0₃0955'
$ws.Range("E28").Value = '  -11.29%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.23%  '
$ws.Range("D30").Value = '6.93'
$ws.Range("E30").Value = '  -5.31%  '
$ws.Range("D31").Value = '0.997'
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("E32").Value = '  -9.10%  '
$ws.Range("D33").Value = '6.98'
$ws.Range("E33").Value = '  -9.57%  '
$ws.Range("D34").Value = '21.69'
$ws.Range("E34").Value = '  -4.82%  '
$ws.Range("D35").Value = '1.21'
$ws.Range("E35").Value = '  -6.98%  '
$ws.Range("D36").Value = '4.96'
$ws.Range("E36").Value = '  -6.72%  '
$ws.Range("D37").Value = '160.13'
$ws.Range("E37").Value = '  -5.01%  '
$ws.Range("D38").Value = '1.43'
$ws.Range("E38").Value = '  -7.69%  '
$ws.Range("D39").Value = '6.27'
$ws.Range("E39").Value = '  -7.90%  '
$ws.Range("E40").Value = '  -7.69%  '
$ws.Range("D41").Value = '0.0704'
$ws.Range("E41").Value = '  -6.09%  '
$ws.Range("D42").Value = '3.204.34'
$ws.Range("E42").Value = '  -5.34%  '
$ws.Range("D43").Value = '40.31'
$ws.Range("E43").Value = '  -4.86%  '
$ws.Range("D44").Value = '0.705'
$ws.Range("E44").Value = '  -7.11%  '
$ws.Range("B45").Value = 'ONDO'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D45").Value = '1.09'
$ws.Range("E45").Value = '  -3.81%  '
$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.00'
$ws.Range("E46").Value = '  -6.83%  '
$ws.Range("E47").Value = '  -0.04%  '
$ws.Range("D48").Value = '1.48'
$ws.Range("E48").Value = '  -8.18%  '
$ws.Range("D49").Value = '2.277.00'
$ws.Range("E49").Value = '  -7.44%  '
$ws.Range("D50").Value = '6.22'
$ws.Range("E50").Value = '  -6.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.70'
$ws.Range("E51").Value = '  -6.79%  '
